$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-6 (columns B-H; column A "sr No." is untouched by the edit)
# B = Delivery No. (new), C = Source code, D = Destination Code,
# E = Transporter Code, F = Vehicle Number, G = Start Date, H = End date
$data = @(
    @{ Row=2; B=10; C="S001"; D="D001"; E="T001"; F="MH04DE1433"; G="2020-04-15"; H="2020-04-15" },
    @{ Row=3; B=11; C="S006"; D="D005"; E="T002"; F="mh05fr3434"; G="2020-04-16"; H="2020-04-16" },
    @{ Row=4; B=12; C="S001"; D="D002"; E="T001"; F="mh462390";   G="2020-02-11"; H="2020-02-11" },
    @{ Row=5; B=13; C="S001"; D="D002"; E="T001"; F="mh083567";   G="2020-04-18"; H="2020-04-18" },
    @{ Row=6; B=14; C="S001"; D="D002"; E="T001"; F="mh083567";   G="2020-02-18"; H="2020-02-18" }
)

foreach ($rowData in $data) {
    $r = $rowData.Row

    # Numeric column
    $ws.Cells.Item($r, 2).Value = $rowData.B

    # Plain text columns (no risk of auto type conversion)
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F

    # Date-like text columns: use Formula + Copy + PasteSpecial(values) so Excel
    # stores them as plain text shared-strings instead of auto-converting to date
    # serial numbers (which would also introduce new number-format styles).
    $ws.Cells.Item($r, 7).Formula = "=""" + $rowData.G + """"
    $ws.Cells.Item($r, 8).Formula = "=""" + $rowData.H + """"
}

$dataRange = $ws.Range("G2:H6")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
